$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws1 = $wb.Worksheets.Item("Rushing")

# Row 2: K.Murray
$ws1.Range("C2").Value = 36
$ws1.Range("D2").Value = 29
$ws1.Range("E2").Value = 16

# Row 5: C.Edmonds
$ws1.Range("C5").Value = 89
$ws1.Range("D5").Value = 32
$ws1.Range("E5").Value = 18

# Row 6: J.Conner
$ws1.Range("C6").Value = 99
$ws1.Range("D6").Value = 52
$ws1.Range("E6").Value = 31
$ws1.Range("F6").Value = 38

# Row 8: E.Benjamin
$ws1.Range("C8").Value = 11
$ws1.Range("D8").Value = 6
$ws1.Range("E8").Value = 4
$ws1.Range("F8").Value = 3

# --- Receiving sheet ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 2: C.Edmonds
$ws2.Range("E2").Value = 3

# Row 3: J.Conner
$ws2.Range("C3").Value = 33
$ws2.Range("D3").Value = 30
$ws2.Range("E3").Value = 4
$ws2.Range("F3").Value = 3
$ws2.Range("G3").Value = 4
$ws2.Range("H3").Value = 3

# Row 4: J.Ward
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 1

# Row 5: E.Benjamin
$ws2.Range("C5").Value = 9
$ws2.Range("D5").Value = 8
$ws2.Range("G5").Value = 1
$ws2.Range("H5").Value = 1

# Row 6: G.Dortch
$ws2.Range("C6").Value = 3
$ws2.Range("D6").Value = 3

# Row 7: A.Green
$ws2.Range("C7").Value = 68
$ws2.Range("D7").Value = 38
$ws2.Range("E7").Value = 28
$ws2.Range("G7").Value = 17
$ws2.Range("H7").Value = 5

# Row 8: C.Kirk
$ws2.Range("C8").Value = 74
$ws2.Range("D8").Value = 58
$ws2.Range("E8").Value = 30
$ws2.Range("F8").Value = 16
$ws2.Range("G8").Value = 9

# Row 9: R.Moore
$ws2.Range("C9").Value = 60
$ws2.Range("D9").Value = 51
$ws2.Range("E9").Value = 6

# Row 10: A.Wesley
$ws2.Range("C10").Value = 21
$ws2.Range("D10").Value = 14
$ws2.Range("E10").Value = 9

# Row 12: M.Williams
$ws2.Range("C12").Value = 5
$ws2.Range("D12").Value = 2

# Row 14: Z.Ertz
$ws2.Range("C14").Value = 82
$ws2.Range("D14").Value = 66
$ws2.Range("E14").Value = 13
$ws2.Range("G14").Value = 13
$ws2.Range("H14").Value = 5
